$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddEmp")

# Clear the old data first (A1:D1 previously held prasad/prasad/123456/123456)
$ws.Range("A1:D1").ClearContents()

# New header/data values
$ws.Range("A1").Value = "VLSI"
$ws.Range("B1").Value = "Chip Designer"
$ws.Range("A2").Value = "Data"
$ws.Range("B2").Value = "DEO"

# Column B custom width
$ws.Columns.Item(2).ColumnWidth = 18.42578125

# New selection on the next empty row
$ws.Range("A3").Select()

